$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 46.17354133333333
$ws.Range("H2").Value = 138.520624
$ws.Range("I2").Value = 0.1154336358852217
$ws.Range("J2").Value = 0.1189208138601986
$ws.Range("M2").Value = 5.343111666666666
$ws.Range("N2").Value = 16.029335
$ws.Range("O2").Value = 0.007536583045055595
$ws.Range("P2").Value = 0.007566622551600167
$ws.Range("Q2").Value = 246.7103873894489
$ws.Range("R2").Value = 2220.39348650504
$ws.Range("S2").Value = 0.0008699751830416832
$ws.Range("T2").Value = 0.0008998289120092245
$ws.Range("G3").Value = 46.17354133333333
$ws.Range("H3").Value = 138.520624
$ws.Range("I3").Value = 0.1154336358852217
$ws.Range("J3").Value = 0.1189208138601986
$ws.Range("O3").Value = 0.03519876756380422
$ws.Range("P3").Value = 0.03533906371688481
$ws.Range("Q3").Value = 1152.23325071624
$ws.Range("R3").Value = 10370.09925644616
$ws.Range("S3").Value = 0.004063121718568729
$ws.Range("T3").Value = 0.004202550218269357
$ws.Range("G4").Value = 46.17354133333333
$ws.Range("H4").Value = 138.520624
$ws.Range("I4").Value = 0.1154336358852217
$ws.Range("J4").Value = 0.1189208138601986
$ws.Range("M4").Value = 312.3302103333334
$ws.Range("N4").Value = 936.990631
$ws.Range("O4").Value = 0.4405490123558179
$ws.Range("P4").Value = 0.4423049639403425
$ws.Range("Q4").Value = 14421.39187647486
$ws.Range("R4").Value = 129792.5268882737
$ws.Range("S4").Value = 0.05085417428187553
$ws.Range("T4").Value = 0.05259926628619133
$ws.Range("G5").Value = 46.17354133333333
$ws.Range("H5").Value = 138.520624
$ws.Range("I5").Value = 0.1154336358852217
$ws.Range("J5").Value = 0.1189208138601986
$ws.Range("M5").Value = 8.4436795
$ws.Range("N5").Value = 16.887359
$ws.Range("O5").Value = 0.01191000595300744
$ws.Range("P5").Value = 0.007971651440709677
$ws.Range("Q5").Value = 389.8745843986694
$ws.Range("R5").Value = 2339.247506392016
$ws.Range("S5").Value = 0.001374815290570284
$ws.Range("T5").Value = 0.0009479952771390195
$ws.Range("G6").Value = 46.17354133333333
$ws.Range("H6").Value = 138.520624
$ws.Range("I6").Value = 0.1154336358852217
$ws.Range("J6").Value = 0.1189208138601986
$ws.Range("M6").Value = 357.885376
$ws.Range("N6").Value = 1073.656128
$ws.Range("O6").Value = 0.504805631082315
$ws.Range("P6").Value = 0.5068176983504629
$ws.Range("Q6").Value = 16524.83520133154
$ws.Range("R6").Value = 148723.5168119839
$ws.Range("S6").Value = 0.05827154941116552
$ws.Range("T6").Value = 0.06027117316658968
$ws.Range("I7").Value = 0.2100327918507284
$ws.Range("J7").Value = 0.2163777511873036
$ws.Range("M7").Value = 5.343111666666666
$ws.Range("N7").Value = 16.029335
$ws.Range("O7").Value = 0.007536583045055595
$ws.Range("P7").Value = 0.007566622551600167
$ws.Range("Q7").Value = 448.892309807375
$ws.Range("R7").Value = 4040.030788266375
$ws.Range("S7").Value = 0.00158292957796789
$ws.Range("T7").Value = 0.001637248771798381
$ws.Range("I8").Value = 0.2100327918507284
$ws.Range("J8").Value = 0.2163777511873036
$ws.Range("O8").Value = 0.03519876756380422
$ws.Range("P8").Value = 0.03533906371688481
$ws.Range("S8").Value = 0.007392895421130661
$ws.Range("T8").Value = 0.00764658713612437
$ws.Range("I9").Value = 0.2100327918507284
$ws.Range("J9").Value = 0.2163777511873036
$ws.Range("M9").Value = 312.3302103333334
$ws.Range("N9").Value = 936.990631
$ws.Range("O9").Value = 0.4405490123558179
$ws.Range("P9").Value = 0.4423049639403425
$ws.Range("Q9").Value = 26239.88385154218
$ws.Range("R9").Value = 236158.9546638796
$ws.Range("S9").Value = 0.09252973901217346
$ws.Range("T9").Value = 0.09570495343639272
$ws.Range("I10").Value = 0.2100327918507284
$ws.Range("J10").Value = 0.2163777511873036
$ws.Range("M10").Value = 8.4436795
$ws.Range("N10").Value = 16.887359
$ws.Range("O10").Value = 0.01191000595300744
$ws.Range("P10").Value = 0.007971651440709677
$ws.Range("Q10").Value = 709.3811678453626
$ws.Range("R10").Value = 4256.287007072176
$ws.Range("S10").Value = 0.002501491801268947
$ws.Range("T10").Value = 0.001724888011989789
$ws.Range("I11").Value = 0.2100327918507284
$ws.Range("J11").Value = 0.2163777511873036
$ws.Range("M11").Value = 357.885376
$ws.Range("N11").Value = 1073.656128
$ws.Range("O11").Value = 0.504805631082315
$ws.Range("P11").Value = 0.5068176983504629
$ws.Range("Q11").Value = 30067.1225123664
$ws.Range("R11").Value = 270604.1026112976
$ws.Range("S11").Value = 0.1060257360381874
$ws.Range("T11").Value = 0.1096640738309983
$ws.Range("G12").Value = 128.0910926666667
$ws.Range("H12").Value = 384.273278
$ws.Range("I12").Value = 0.3202271284388135
$ws.Range("J12").Value = 0.3299009897940278
$ws.Range("M12").Value = 5.343111666666666
$ws.Range("N12").Value = 16.029335
$ws.Range("O12").Value = 0.007536583045055595
$ws.Range("P12").Value = 0.007566622551600167
$ws.Range("Q12").Value = 684.4050116233477
$ws.Range("R12").Value = 6159.64510461013
$ws.Range("S12").Value = 0.002413418346758802
$ws.Range("T12").Value = 0.002496236269170707
$ws.Range("G13").Value = 128.0910926666667
$ws.Range("H13").Value = 384.273278
$ws.Range("I13").Value = 0.3202271284388135
$ws.Range("J13").Value = 0.3299009897940278
$ws.Range("O13").Value = 0.03519876756380422
$ws.Range("P13").Value = 0.03533906371688481
$ws.Range("Q13").Value = 3196.43700329653
$ws.Range("R13").Value = 28767.93302966877
$ws.Range("S13").Value = 0.01127160026154228
$ws.Range("T13").Value = 0.01165839209859451
$ws.Range("G14").Value = 128.0910926666667
$ws.Range("H14").Value = 384.273278
$ws.Range("I14").Value = 0.3202271284388135
$ws.Range("J14").Value = 0.3299009897940278
$ws.Range("M14").Value = 312.3302103333334
$ws.Range("N14").Value = 936.990631
$ws.Range("O14").Value = 0.4405490123558179
$ws.Range("P14").Value = 0.4423049639403425
$ws.Range("Q14").Value = 40006.71791440649
$ws.Range("R14").Value = 360060.4612296585
$ws.Range("S14").Value = 0.1410757451632589
$ws.Range("T14").Value = 0.1459168453947308
$ws.Range("G15").Value = 128.0910926666667
$ws.Range("H15").Value = 384.273278
$ws.Range("I15").Value = 0.3202271284388135
$ws.Range("J15").Value = 0.3299009897940278
$ws.Range("M15").Value = 8.4436795
$ws.Range("N15").Value = 16.887359
$ws.Range("O15").Value = 0.01191000595300744
$ws.Range("P15").Value = 0.007971651440709677
$ws.Range("Q15").Value = 1081.560133282134
$ws.Range("R15").Value = 6489.360799692802
$ws.Range("S15").Value = 0.003813907006020747
$ws.Range("T15").Value = 0.00262985570058311
$ws.Range("G16").Value = 128.0910926666667
$ws.Range("H16").Value = 384.273278
$ws.Range("I16").Value = 0.3202271284388135
$ws.Range("J16").Value = 0.3299009897940278
$ws.Range("M16").Value = 357.885376
$ws.Range("N16").Value = 1073.656128
$ws.Range("O16").Value = 0.504805631082315
$ws.Range("P16").Value = 0.5068176983504629
$ws.Range("Q16").Value = 45841.92886126084
$ws.Range("R16").Value = 412577.3597513476
$ws.Range("S16").Value = 0.1616524576612328
$ws.Range("T16").Value = 0.1671996603309487
$ws.Range("G17").Value = 35.18830149999999
$ws.Range("H17").Value = 70.37660299999999
$ws.Range("I17").Value = 0.08797058803540478
$ws.Range("J17").Value = 0.06041874966919073
$ws.Range("M17").Value = 5.343111666666666
$ws.Range("N17").Value = 16.029335
$ws.Range("O17").Value = 0.007536583045055595
$ws.Range("P17").Value = 0.007566622551600167
$ws.Range("Q17").Value = 188.0150242748341
$ws.Range("R17").Value = 1128.090145649005
$ws.Range("S17").Value = 0.0006629976422512022
$ws.Range("T17").Value = 0.0004571658737863837
$ws.Range("G18").Value = 35.18830149999999
$ws.Range("H18").Value = 70.37660299999999
$ws.Range("I18").Value = 0.08797058803540478
$ws.Range("J18").Value = 0.06041874966919073
$ws.Range("O18").Value = 0.03519876756380422
$ws.Range("P18").Value = 0.03533906371688481
$ws.Range("Q18").Value = 878.1031268931073
$ws.Range("R18").Value = 5268.618761358643
$ws.Range("S18").Value = 0.003096456280709389
$ws.Range("T18").Value = 0.002135142044254044
$ws.Range("G19").Value = 35.18830149999999
$ws.Range("H19").Value = 70.37660299999999
$ws.Range("I19").Value = 0.08797058803540478
$ws.Range("J19").Value = 0.06041874966919073
$ws.Range("M19").Value = 312.3302103333334
$ws.Range("N19").Value = 936.990631
$ws.Range("O19").Value = 0.4405490123558179
$ws.Range("P19").Value = 0.4423049639403425
$ws.Range("Q19").Value = 10990.36960876775
$ws.Range("R19").Value = 65942.21765260649
$ws.Range("S19").Value = 0.03875535567535811
$ws.Range("T19").Value = 0.02672351289375199
$ws.Range("G20").Value = 35.18830149999999
$ws.Range("H20").Value = 70.37660299999999
$ws.Range("I20").Value = 0.08797058803540478
$ws.Range("J20").Value = 0.06041874966919073
$ws.Range("M20").Value = 8.4436795
$ws.Range("N20").Value = 16.887359
$ws.Range("O20").Value = 0.01191000595300744
$ws.Range("P20").Value = 0.007971651440709677
$ws.Range("Q20").Value = 297.1187400153692
$ws.Range("R20").Value = 1188.474960061477
$ws.Range("S20").Value = 0.001047730227191236
$ws.Range("T20").Value = 0.0004816372128462816
$ws.Range("G21").Value = 35.18830149999999
$ws.Range("H21").Value = 70.37660299999999
$ws.Range("I21").Value = 0.08797058803540478
$ws.Range("J21").Value = 0.06041874966919073
$ws.Range("M21").Value = 357.885376
$ws.Range("N21").Value = 1073.656128
$ws.Range("O21").Value = 0.504805631082315
$ws.Range("P21").Value = 0.5068176983504629
$ws.Range("Q21").Value = 12593.37851312886
$ws.Range("R21").Value = 75560.27107877318
$ws.Range("S21").Value = 0.04440804820989486
$ws.Range("T21").Value = 0.03062129164455204
$ws.Range("G22").Value = 106.534543
$ws.Range("H22").Value = 319.603629
$ws.Range("I22").Value = 0.2663358557898317
$ws.Range("J22").Value = 0.2743816954892795
$ws.Range("M22").Value = 5.343111666666666
$ws.Range("N22").Value = 16.029335
$ws.Range("O22").Value = 0.007536583045055595
$ws.Range("P22").Value = 0.007566622551600167
$ws.Range("Q22").Value = 569.2259596063017
$ws.Range("R22").Value = 5123.033636456715
$ws.Range("S22").Value = 0.002007262295036017
$ws.Range("T22").Value = 0.002076142724835472
$ws.Range("G23").Value = 106.534543
$ws.Range("H23").Value = 319.603629
$ws.Range("I23").Value = 0.2663358557898317
$ws.Range("J23").Value = 0.2743816954892795
$ws.Range("O23").Value = 0.03519876756380422
$ws.Range("P23").Value = 0.03533906371688481
$ws.Range("Q23").Value = 2658.506132511915
$ws.Range("R23").Value = 23926.55519260724
$ws.Range("S23").Value = 0.009374693881853163
$ws.Range("T23").Value = 0.009696392219642533
$ws.Range("G24").Value = 106.534543
$ws.Range("H24").Value = 319.603629
$ws.Range("I24").Value = 0.2663358557898317
$ws.Range("J24").Value = 0.2743816954892795
$ws.Range("M24").Value = 312.3302103333334
$ws.Range("N24").Value = 936.990631
$ws.Range("O24").Value = 0.4405490123558179
$ws.Range("P24").Value = 0.4423049639403425
$ws.Range("Q24").Value = 33273.95622295554
$ws.Range("R24").Value = 299465.6060065999
$ws.Range("S24").Value = 0.1173339982231519
$ws.Range("T24").Value = 0.1213603859292758
$ws.Range("G25").Value = 106.534543
$ws.Range("H25").Value = 319.603629
$ws.Range("I25").Value = 0.2663358557898317
$ws.Range("J25").Value = 0.2743816954892795
$ws.Range("M25").Value = 8.4436795
$ws.Range("N25").Value = 16.887359
$ws.Range("O25").Value = 0.01191000595300744
$ws.Range("P25").Value = 0.007971651440709677
$ws.Range("Q25").Value = 899.5435367709684
$ws.Range("R25").Value = 5397.261220625811
$ws.Range("S25").Value = 0.003172061627956226
$ws.Range("T25").Value = 0.002187275238151479
$ws.Range("G26").Value = 106.534543
$ws.Range("H26").Value = 319.603629
$ws.Range("I26").Value = 0.2663358557898317
$ws.Range("J26").Value = 0.2743816954892795
$ws.Range("M26").Value = 357.885376
$ws.Range("N26").Value = 1073.656128
$ws.Range("O26").Value = 0.504805631082315
$ws.Range("P26").Value = 0.5068176983504629
$ws.Range("Q26").Value = 38127.15497854317
$ws.Range("R26").Value = 343144.3948068885
$ws.Range("S26").Value = 0.1344478397618344
$ws.Range("T26").Value = 0.1390614993773742
